$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "RT @TribulationThe: TREASON by the Biden Crime Family!`n`"Please Have Keys Made`": Joe Biden Was Chinese Financier's `"Office Mate`" According…"
$ws.Range("C2").Value = [double]"1.339313651711357e+18"
$ws.Range("D2").Value = [double]"139"
$ws.Range("E2").Value = [double]"44181.87305555555"
$ws.Range("F2").Value = "en"
$ws.Range("G2").Value = [double]"22442937"
$ws.Range("H2").Value = "TheTybeeTimes"
$ws.Range("I2").Value = [double]"2230"
$ws.Range("J2").Value = "Tybee Island, GA"
$ws.Range("K2").Value = "RT @TribulationThe: TREASON by the Biden Crime Family!`"Please Have Keys Made`": Joe Biden Was Chinese Financier's `"Office Mate`" According… "
$ws.Range("L2").Value = " TREASON by the Biden Crime Family!  `"Please Have Keys Made`": Joe Biden Was Chinese Financier's `"Office Mate`" According…"
$ws.Range("M2").Value = " TREASON by the Biden Crime Family!`"Please Have Keys Made`": Joe Biden Was Chinese Financier's `"Office Mate`" According… "
$ws.Range("N2").Value = "RT @TribulationThe: TREASON by the Biden Crime Family!`"Please Have Keys Made`": Joe Biden Was Chinese Financier's `"Office Mate`" According… "
# Row 3
$ws.Range("B3").Value = "Sen. Mitch McConnell: `"Under President #Trump's command, our forces took terrorist leaders like Baghdadi &amp;  Qasem S… https://t.co/LRH8KmSy2k"
$ws.Range("C3").Value = [double]"1.33931361351388e+18"
$ws.Range("D3").Value = [double]"144"
$ws.Range("E3").Value = [double]"44181.87295138889"
$ws.Range("F3").Value = "en"
$ws.Range("G3").Value = [double]"2343557479"
$ws.Range("H3").Value = "MostafaMe4"
$ws.Range("I3").Value = [double]"7816"
$ws.Range("J3").Value = "Europe"
$ws.Range("K3").Value = "Sen. Mitch McConnell: `"Under President #Trump's command, our forces took terrorist leaders like Baghdadi &amp;  Qasem S… https://t.co/LRH8KmSy2k "
$ws.Range("L3").Value = "Sen. Mitch McConnell: `"Under President #Trump's command, our forces took terrorist leaders like Baghdadi &amp;  Qasem S… "
$ws.Range("M3").Value = "Sen. Mitch McConnell: `"Under President #Trump's command, our forces took terrorist leaders like Baghdadi &amp;  Qasem S…  "
$ws.Range("N3").Value = "Sen. Mitch McConnell: `"Under President #Trump's command, our forces took terrorist leaders like Baghdadi &amp;  Qasem S…  "
# Row 4
$ws.Range("B4").Value = "RT @yoksig: @SusanLynch22 @SuzieBird4 @Missin_Florida @skewermann @P4boxers @HeathenResister @openpodbaydoor_ @DearAuntCrabby @BlogChurchWo…"
$ws.Range("C4").Value = [double]"1.339313610762445e+18"
$ws.Range("D4").Value = [double]"140"
$ws.Range("E4").Value = [double]"44181.87293981481"
$ws.Range("F4").Value = "en"
$ws.Range("G4").Value = [double]"247054633"
$ws.Range("H4").Value = "jkf3500"
$ws.Range("I4").Value = [double]"25089"
$ws.Range("J4").Value = "Texas"
$ws.Range("K4").Value = "RT @yoksig: @SusanLynch22 @SuzieBird4 @Missin_Florida @skewermann @P4boxers @HeathenResister @openpodbaydoor_ @DearAuntCrabby @BlogChurchWo… "
$ws.Range("L4").Value = "         …"
$ws.Range("M4").Value = "         … "
$ws.Range("N4").Value = "RT @yoksig: @SusanLynch22 @SuzieBird4 @Missin_Florida @skewermann @P4boxers @HeathenResister @openpodbaydoor_ @DearAuntCrabby @BlogChurchWo… "
# Row 5
$ws.Range("B5").Value = "Nenhum dos repórteres na sala perguntou-lhe sobre o espião chinês no escritório de Eric Swalwell, então ela mesma f… https://t.co/ijmucNb4Ng"
$ws.Range("C5").Value = [double]"1.339313569838588e+18"
$ws.Range("D5").Value = [double]"140"
$ws.Range("E5").Value = [double]"44181.87282407407"
$ws.Range("F5").Value = "pt"
$ws.Range("G5").Value = [double]"128376943"
$ws.Range("H5").Value = "German_Emanuel"
$ws.Range("I5").Value = [double]"61"
# J5: left empty (no change)
$ws.Range("K5").Value = "None of the reporters in the room asked her about the Chinese spy at Eric Swalwell's office, so she herself… https://t.co/ijmucNb4Ng "
$ws.Range("L5").Value = "Nenhum dos repórteres na sala perguntou-lhe sobre o espião chinês no escritório de Eric Swalwell, então ela mesma f… "
$ws.Range("M5").Value = "None of the reporters in the room asked her about the Chinese spy at Eric Swalwell's office, so she herself…  "
$ws.Range("N5").Value = "None of the reporters in the room asked her about the Chinese spy at Eric Swalwell's office, so she herself…  "
# Row 6
$ws.Range("B6").Value = "Apparently Trump thought he was running against himself not against Biden. No, beating your own score doesn’t autom… https://t.co/89zD1fedtN"
$ws.Range("C6").Value = [double]"1.339313511420183e+18"
$ws.Range("D6").Value = [double]"140"
$ws.Range("E6").Value = [double]"44181.87267361111"
$ws.Range("F6").Value = "en"
$ws.Range("G6").Value = [double]"1.017492973523698e+18"
$ws.Range("H6").Value = "RachelWilder_"
$ws.Range("I6").Value = [double]"229"
$ws.Range("J6").Value = "🌎"
$ws.Range("K6").Value = "Apparently Trump thought he was running against himself not against Biden. No, beating your own score doesn’t autom… https://t.co/89zD1fedtN "
$ws.Range("L6").Value = "Apparently Trump thought he was running against himself not against Biden. No, beating your own score doesn’t autom… "
$ws.Range("M6").Value = "Apparently Trump thought he was running against himself not against Biden. No, beating your own score doesn’t autom…  "
$ws.Range("N6").Value = "Apparently Trump thought he was running against himself not against Biden. No, beating your own score doesn’t autom…  "
# Row 7
$ws.Range("B7").Value = "RT @AllenLEllison: .@StephenKing even you couldn’t have written a book more horrific than what #Trump has done to this country. #Rubio had…"
$ws.Range("C7").Value = [double]"1.339313501899088e+18"
$ws.Range("D7").Value = [double]"139"
$ws.Range("E7").Value = [double]"44181.87263888889"
$ws.Range("F7").Value = "en"
$ws.Range("G7").Value = [double]"2508838914"
$ws.Range("H7").Value = "StuTheJanitor"
$ws.Range("I7").Value = [double]"4738"
$ws.Range("J7").Value = "Portland, OR"
$ws.Range("K7").Value = "RT @AllenLEllison: .@StephenKing even you couldn’t have written a book more horrific than what #Trump has done to this country. #Rubio had… "
$ws.Range("L7").Value = " . even you couldn’t have written a book more horrific than what #Trump has done to this country. #Rubio had…"
$ws.Range("M7").Value = " . even you couldn’t have written a book more horrific than what #Trump has done to this country. #Rubio had… "
$ws.Range("N7").Value = "RT @AllenLEllison: .@StephenKing even you couldn’t have written a book more horrific than what #Trump has done to this country. #Rubio had… "
# Row 8
$ws.Range("B8").Value = "RT @PdS1748: @Butterf70713546 @Dragonmaster969 @Boduoghnat @telegraaf Wat een geweldig artikel over de handelingen van #Trump. #TRUMP2020To…"
$ws.Range("C8").Value = [double]"1.339313472287498e+18"
$ws.Range("D8").Value = [double]"140"
$ws.Range("E8").Value = [double]"44181.87255787037"
$ws.Range("F8").Value = "nl"
$ws.Range("G8").Value = [double]"1.276310374669783e+18"
$ws.Range("H8").Value = "Rudy84529712"
$ws.Range("I8").Value = [double]"619"
# J8: left empty (no change)
$ws.Range("K8").Value = "RT @ PdS1748: @ Butterf70713546 @ Dragonmaster969 @Boduoghnat @telegraaf What a great article about the actions of #Trump. # TRUMP2020To ... "
$ws.Range("L8").Value = "     Wat een geweldig artikel over de handelingen van #Trump. #TRUMP2020To…"
$ws.Range("M8").Value = "RT  PdS1748:  Butterf70713546  Dragonmaster969 Boduoghnat telegraaf What a great article about the actions of #Trump. # TRUMP2020To ... "
$ws.Range("N8").Value = "RT @ PdS1748: @ Butterf70713546 @ Dragonmaster969 @Boduoghnat @telegraaf What a great article about the actions of #Trump. # TRUMP2020To ... "
# Row 9
$ws.Range("B9").Value = "#MAGA #Trump #Republicans #GOP #MitchMcConnell #Canada #StimulusCheckNOW #Stimuluscheck #JoseBiden #FoxNews #COVID19 https://t.co/gP6k3ATA2q"
$ws.Range("C9").Value = [double]"1.339313444781093e+18"
$ws.Range("D9").Value = [double]"140"
$ws.Range("E9").Value = [double]"44181.87248842593"
$ws.Range("F9").Value = "und"
$ws.Range("G9").Value = [double]"1.337280942184878e+18"
$ws.Range("H9").Value = "cynthia26942734"
$ws.Range("I9").Value = [double]"16"
# J9: left empty (no change)
$ws.Range("K9").Value = "#MAGA #Trump #Republicans #GOP #MitchMcConnell #Canada #StimulusCheckNOW #Stimuluscheck #JoseBiden #FoxNews #COVID19 https://t.co/gP6k3ATA2q "
$ws.Range("L9").Value = "#MAGA #Trump #Republicans #GOP #MitchMcConnell #Canada #StimulusCheckNOW #Stimuluscheck #JoseBiden #FoxNews #COVID19 "
$ws.Range("M9").Value = "#MAGA #Trump #Republicans #GOP #MitchMcConnell #Canada #StimulusCheckNOW #Stimuluscheck #JoseBiden #FoxNews #COVID19  "
$ws.Range("N9").Value = "#MAGA #Trump #Republicans #GOP #MitchMcConnell #Canada #StimulusCheckNOW #Stimuluscheck #JoseBiden #FoxNews #COVID19  "
# Row 10
$ws.Range("B10").Value = "#Trump literally gave a foreign ENEMY state that has nuclear weapons pointed at us backdoor access into ALL governm… https://t.co/QGT9RpEjfK"
$ws.Range("C10").Value = [double]"1.339313434664505e+18"
$ws.Range("D10").Value = [double]"140"
$ws.Range("E10").Value = [double]"44181.87245370371"
$ws.Range("F10").Value = "en"
$ws.Range("G10").Value = [double]"9.918115387125146e+17"
$ws.Range("H10").Value = "_Anonym0us_FL_"
$ws.Range("I10").Value = [double]"138"
$ws.Range("J10").Value = "Florida, USA"
$ws.Range("K10").Value = "#Trump literally gave a foreign ENEMY state that has nuclear weapons pointed at us backdoor access into ALL governm… https://t.co/QGT9RpEjfK "
$ws.Range("L10").Value = "#Trump literally gave a foreign ENEMY state that has nuclear weapons pointed at us backdoor access into ALL governm… "
$ws.Range("M10").Value = "#Trump literally gave a foreign ENEMY state that has nuclear weapons pointed at us backdoor access into ALL governm…  "
$ws.Range("N10").Value = "#Trump literally gave a foreign ENEMY state that has nuclear weapons pointed at us backdoor access into ALL governm…  "
# Row 11
$ws.Range("B11").Value = "President #Trump did in fact win the Election `nفاکس نیوز هم اقرار به بُرد آقای ترامپ در انتخابات شد.… https://t.co/OdA7VTt2XE"
$ws.Range("C11").Value = [double]"1.339313427144012e+18"
$ws.Range("D11").Value = [double]"125"
$ws.Range("E11").Value = [double]"44181.87243055556"
$ws.Range("F11").Value = "und"
$ws.Range("G11").Value = [double]"9.723674191161467e+17"
$ws.Range("H11").Value = "Ali49781135"
$ws.Range("I11").Value = [double]"1094"
$ws.Range("J11").Value = "Iran,Tehran"
$ws.Range("K11").Value = "President #Trump did in fact win the Election Fox News also acknowledged Mr. Trump's victory in the election.… Https://t.co/OdA7VTt2XE "
$ws.Range("L11").Value = "President #Trump did in fact win the Election  فاکس نیوز هم اقرار به بُرد آقای ترامپ در انتخابات شد.… "
$ws.Range("M11").Value = "President #Trump did in fact win the Election Fox News also acknowledged Mr. Trump's victory in the election.… Https://t.co/OdA7VTt2XE "
$ws.Range("N11").Value = "President #Trump did in fact win the Election Fox News also acknowledged Mr. Trump's victory in the election.… Https://t.co/OdA7VTt2XE "
